{"js": "// Regenerate the 25 division-problem answers in the single table (5 data\n// rows x 5 columns) with a new batch of problems, cell by cell, preserving\n// each run's existing formatting (font/size/paragraph alignment).\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"items\");\nawait context.sync();\n\nconst updates = [\n  { row: 0, col: 0, oldText: \"35\u00f73=11, 2\", newText: \"25\u00f72=12, 1\" },\n  { row: 0, col: 1, oldText: \"75\u00f72=37, 1\", newText: \"10\u00f74=2, 2\" },\n  { row: 0, col: 2, oldText: \"72\u00f74=18, 0\", newText: \"42\u00f78=5, 2\" },\n  { row: 0, col: 3, oldText: \"51\u00f74=12, 3\", newText: \"93\u00f77=13, 2\" },\n  { row: 0, col: 4, oldText: \"83\u00f72=41, 1\", newText: \"39\u00f73=13, 0\" },\n  { row: 4, col: 0, oldText: \"90\u00f77=12, 6\", newText: \"98\u00f77=14, 0\" },\n  { row: 4, col: 1, oldText: \"48\u00f74=12, 0\", newText: \"58\u00f78=7, 2\" },\n  { row: 4, col: 2, oldText: \"89\u00f78=11, 1\", newText: \"15\u00f73=5, 0\" },\n  { row: 4, col: 3, oldText: \"40\u00f77=5, 5\", newText: \"95\u00f78=11, 7\" },\n  { row: 4, col: 4, oldText: \"31\u00f76=5, 1\", newText: \"48\u00f72=24, 0\" },\n  { row: 8, col: 0, oldText: \"15\u00f75=3, 0\", newText: \"27\u00f78=3, 3\" },\n  { row: 8, col: 1, oldText: \"90\u00f78=11, 2\", newText: \"63\u00f73=21, 0\" },\n  { row: 8, col: 2, oldText: \"44\u00f77=6, 2\", newText: \"69\u00f72=34, 1\" },\n  { row: 8, col: 3, oldText: \"48\u00f74=12, 0\", newText: \"70\u00f76=11, 4\" },\n  { row: 8, col: 4, oldText: \"76\u00f74=19, 0\", newText: \"91\u00f79=10, 1\" },\n  { row: 12, col: 0, oldText: \"24\u00f74=6, 0\", newText: \"15\u00f72=7, 1\" },\n  { row: 12, col: 1, oldText: \"44\u00f73=14, 2\", newText: \"21\u00f75=4, 1\" },\n  { row: 12, col: 2, oldText: \"88\u00f72=44, 0\", newText: \"56\u00f72=28, 0\" },\n  { row: 12, col: 3, oldText: \"93\u00f73=31, 0\", newText: \"93\u00f72=46, 1\" },\n  { row: 12, col: 4, oldText: \"45\u00f76=7, 3\", newText: \"73\u00f74=18, 1\" },\n  { row: 16, col: 0, oldText: \"39\u00f75=7, 4\", newText: \"34\u00f79=3, 7\" },\n  { row: 16, col: 1, oldText: \"72\u00f74=18, 0\", newText: \"34\u00f76=5, 4\" },\n  { row: 16, col: 2, oldText: \"53\u00f77=7, 4\", newText: \"18\u00f76=3, 0\" },\n  { row: 16, col: 3, oldText: \"55\u00f72=27, 1\", newText: \"64\u00f79=7, 1\" },\n  { row: 16, col: 4, oldText: \"75\u00f73=25, 0\", newText: \"36\u00f76=6, 0\" },\n];\n\nfor (const u of updates) {\n  const cell = table.getCell(u.row, u.col);\n  const results = cell.body.search(u.oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Expected text not found in cell (${u.row}, ${u.col}): ${u.oldText}`);\n  }\n  results.items[0].insertText(u.newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the generated division-problem answers in the table to a new\n# random batch (25 cells across 5 data rows x 5 columns).\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$t.Cell(1,1).Range.Text = \"25\u00f72=12, 1\"\n$t.Cell(1,2).Range.Text = \"10\u00f74=2, 2\"\n$t.Cell(1,3).Range.Text = \"42\u00f78=5, 2\"\n$t.Cell(1,4).Range.Text = \"93\u00f77=13, 2\"\n$t.Cell(1,5).Range.Text = \"39\u00f73=13, 0\"\n$t.Cell(5,1).Range.Text = \"98\u00f77=14, 0\"\n$t.Cell(5,2).Range.Text = \"58\u00f78=7, 2\"\n$t.Cell(5,3).Range.Text = \"15\u00f73=5, 0\"\n$t.Cell(5,4).Range.Text = \"95\u00f78=11, 7\"\n$t.Cell(5,5).Range.Text = \"48\u00f72=24, 0\"\n$t.Cell(9,1).Range.Text = \"27\u00f78=3, 3\"\n$t.Cell(9,2).Range.Text = \"63\u00f73=21, 0\"\n$t.Cell(9,3).Range.Text = \"69\u00f72=34, 1\"\n$t.Cell(9,4).Range.Text = \"70\u00f76=11, 4\"\n$t.Cell(9,5).Range.Text = \"91\u00f79=10, 1\"\n$t.Cell(13,1).Range.Text = \"15\u00f72=7, 1\"\n$t.Cell(13,2).Range.Text = \"21\u00f75=4, 1\"\n$t.Cell(13,3).Range.Text = \"56\u00f72=28, 0\"\n$t.Cell(13,4).Range.Text = \"93\u00f72=46, 1\"\n$t.Cell(13,5).Range.Text = \"73\u00f74=18, 1\"\n$t.Cell(17,1).Range.Text = \"34\u00f79=3, 7\"\n$t.Cell(17,2).Range.Text = \"34\u00f76=5, 4\"\n$t.Cell(17,3).Range.Text = \"18\u00f76=3, 0\"\n$t.Cell(17,4).Range.Text = \"64\u00f79=7, 1\"\n$t.Cell(17,5).Range.Text = \"36\u00f76=6, 0\"\n"}
